# Update party labels so that coalition/party names carry their full
# descriptive form (accounts for parties whose short code could contain a
# "+" once combined into coalition keys). This mirrors a change to the
# shared-string table: every cell whose text is one of the old party /
# coalition labels gets replaced by the corresponding new label, across all
# worksheets in the workbook.

$wb = $excel.ActiveWorkbook

$map = @{
    "Dem"           = "Dem - Democratic Party (Democratic Party, Dem)"
    "Oth"           = "Other - Other"
    "Rep"           = "Rep - Republican Party (Republican Party, Rep)"
    "Dem+Oth"       = "Dem - Democratic Party (Democratic Party, Dem)+Other - Other"
    "Dem+Rep"       = "Dem - Democratic Party (Democratic Party, Dem)+Rep - Republican Party (Republican Party, Rep)"
    "Oth+Rep"       = "Other - Other+Rep - Republican Party (Republican Party, Rep)"
    "Dem+Oth+Rep"   = "Dem - Democratic Party (Democratic Party, Dem)+Other - Other+Rep - Republican Party (Republican Party, Rep)"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) { continue }

    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value()
            if ($val -ne $null -and $map.ContainsKey([string]$val)) {
                $cell.Value = $map[[string]$val]
            }
        }
    }
}
